# Apply the "database filled with data" update to the Database tables workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Correct the misspelled street name "Eslweide" -> "Elsweide" used by testperson #4 (row 16).
$ws.Range("H16").Value = "Elsweide"

# sensordata rows: the sensor value storage path now points at the sensor data folder
# instead of holding a bare numeric placeholder.
$ws.Range("I21").Value = "/data/sensors/"
$ws.Range("I22").Value = "/data/sensors/"
$ws.Range("I23").Value = "/data/sensors/"
$ws.Range("I24").Value = "/data/sensors/"
$ws.Range("I25").Value = "/data/sensors/"

# Add sensor_unit column (header + units) to the "sensor" table block (rows 13-16).
$ws.Range("M13").Value = "sensor_unit"
$ws.Range("M14").Value = "°C"
$ws.Range("M15").Value = "%"

# Update the active selection to reflect where the author left off working.
$ws.Range("N13").Select()
